$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47; existing rows 47-61 shift down to 48-62.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly price observation.
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C47").Value = "Arica y Parinacota"
$ws.Range("D47").Value = 45233
$ws.Range("E47").Value = 15
$ws.Range("F47").Value = 100112044
$ws.Range("G47").Value = "Perejil"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 650
$ws.Range("K47").Value = 1000
$ws.Range("L47").Value = 1300
$ws.Range("M47").Value = 1138
$ws.Range("N47").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O47").Value = "Región de Arica y Parinacota"
$ws.Range("P47").Value = 569
$ws.Range("Q47").Value = 2
$ws.Range("R47").Value = "Hortaliza"
